# Update product collection data for June 2018.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (sku / shopify ids) ------------------------------------------------
$ws.Range("A2").Value = 764204172474
$ws.Range("C2").Value = 224674119698
$ws.Range("D2").Value = 2008926027794

# --- Row 3 (sku / shopify ids) ------------------------------------------------
$ws.Range("A3").Value = 764204172467
$ws.Range("C3").Value = 224674086930
$ws.Range("D3").Value = 2008925995026

# --- product_title column (accented "Rosé") -----------------------------------
$ws.Range("B2").Value = "Rosé All Day - 3 Items"
$ws.Range("B3").Value = "Rosé All Day - 5 Items"

# --- product_collection column (plain "Rose") ---------------------------------
$ws.Range("E3").Value = "Rose All Day - 5 Items"
$ws.Range("E2").Value = "Rose All Day - 3 Items"

# Move the active selection to A2, matching the saved workbook view state.
$null = $ws.Range("A2").Select()
